# Update "Código Efecto Hall" workbook:
#  - Rename the VL_* headers' unit from (mV) to (V), and relabel column E as
#    B_err_rounded (mT) (shifting the VL_* headers one column to the right,
#    F..J).
#  - Change the VL_err column (J) data values from 0.1 to 0.001 (1E-3),
#    consistent with the new Volt unit.
#  - Leave the sheet positioned/selected on the edited J2:J16 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Header row relabeling (row 1, columns E:J) ---
$ws.Range("E1").Value = "B_err_rounded (mT)"
$ws.Range("F1").Value = "VL_-20mA (V)"
$ws.Range("G1").Value = "VL_-10mA (V)"
$ws.Range("H1").Value = "VL_10mA (V)"
$ws.Range("I1").Value = "VL_20mA (V)"
$ws.Range("J1").Value = "VL_err (V)"

# --- VL_err (column J) data rows: 0.1 -> 0.001 ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.001
}

# --- Leave the view focused on the column that changed ---
$ws.Activate()
$ws.Range("J2:J16").Select()
